$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 20 (shifts rows 20..57 down to 21..58,
# growing the used range from A1:R57 to A1:R58).
$ws.Rows("20").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44665
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 100112043
$ws.Range("G20").Value = "Pepino dulce"
$ws.Range("H20").Value = "Cultivar IV Región"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 22000
$ws.Range("L20").Value = 22000
$ws.Range("M20").Value = 22000
$ws.Range("N20").Value = "$/bandeja 18 kilos"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 1222
$ws.Range("Q20").Value = 18
$ws.Range("R20").Value = "Hortaliza"
